# Appends the "JPA" section (heading + body paragraph about persist/optional/nullable)
# after the existing JWT/Spring-Security troubleshooting notes, preceded by two
# blank paragraphs, matching the target OOXML structure exactly.

$d = $word.ActiveDocument

$insertionPoint = $d.Content
$insertionPoint.Collapse(0)  # wdCollapseEnd

$newContentXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>JPA</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">A la hora de hacer un </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>persist</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> con una clase que tenga varios hijos y estos hijos no tengan restricción de existencia  el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>persist</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> fallara si el valor </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>optional</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">= false o </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>nullable</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">=false o a la hora d </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ela</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> creación en </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bd</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> se puso que la columna no podía ser nula.</w:t></w:r></w:p>'

$null = $insertionPoint.InsertXML($newContentXml)
